$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.015.94"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.466.54"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.42"
$ws.Range("E5").Value = "  -4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.92"
$ws.Range("E6").Value = "  -4.45%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0992"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.40"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.905.37"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "57.922.04"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.04"
$ws.Range("E15").Value = "  -4.99%  "
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.470.30"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.85"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.16"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "318.88"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.76"
$ws.Range("E22").Value = "  -3.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.19"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.408"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.39"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0746"
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("E29").Value = "  -5.07%  "
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.72"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.12"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("E36").Value = "  -9.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.95"
$ws.Range("E37").Value = "  -4.39%  "
$ws.Range("E38").Value = "  -4.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.789"
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.48"
$ws.Range("E40").Value = "  -4.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "274.55"
$ws.Range("E41").Value = "  -5.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.95"
$ws.Range("E42").Value = "  -5.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.88"
$ws.Range("E44").Value = "  -3.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0909"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("E46").Value = "  -3.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0214"
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.03"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.733.91"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.65"
$ws.Range("E51").Value = "  -2.51%  "
